$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    0.1687277178753042,
    0.09280387088913836,
    0.07533787605844625,
    0.07353258278773252,
    0.06522649323474324,
    0.04815114568919754,
    0.04544175536003852,
    0.04520291902481287,
    0.04224598982189771,
    0.0396342613140253,
    0.03277762790072531,
    0.03056428226780419,
    0.02677094241723644,
    0.02241870620326616,
    0.02169156076160749,
    0.02124579739910132,
    0.01933897451219854,
    0.01583185148441133,
    0.01572935491949668,
    0.01444428006316109,
    0.01417160052253909,
    0.01329264413095967,
    0.0088736884549223,
    0.008764423248928378,
    0.008139774372184366,
    0.007582811906233048,
    0.006711591104458504,
    0.005654716147366854,
    0.00397706010843363,
    0.002378693865000421,
    0.002327445582543095,
    0.0008431792887318565,
    0.0001179677445244111,
    0.0000464135388292765
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}
